# Commit message: "update data with resort sheetname"
# The workbook has two sheets: "2022-Q2" (big fund table) and "总计" (summary table).
# The edit re-sorts the sheet tabs so that "总计" comes first, before "2022-Q2".

$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item("总计")
$dataSheet = $wb.Worksheets.Item("2022-Q2")

# Move "总计" so it becomes the first tab, before "2022-Q2".
$summarySheet.Move($dataSheet)

# Keep "2022-Q2" as the active/selected tab after the reorder (re-fetch it
# fresh from the Worksheets collection so the selection is applied to the
# sheet's current, post-move state).
$wb.Worksheets.Item("2022-Q2").Select()
